$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2162.75
$ws.Range("I137").Value = 2220.4
$ws.Range("J137").Value = 2066.6667
$ws.Range("K137").Value = 6661.200000000001
$ws.Range("L137").Value = 6200.000100000001
$ws.Range("M137").Value = -4111.200000000001
$ws.Range("N137").Value = -11300.0001

$ws.Range("H141").Value = 6291.615
$ws.Range("I141").Value = 3657.2
$ws.Range("K141").Value = 10971.6
$ws.Range("M141").Value = -5791.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9805476
$ws.Range("I61").Value = 9805476
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 9805476
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -9805264
$ws.Range("N61").ClearContents()

$ws.Range("H63").Value = 2448.795
$ws.Range("I63").Value = 2173
$ws.Range("J63").Value = 2557.1428
$ws.Range("K63").Value = 2173
$ws.Range("L63").Value = 2557.1428
$ws.Range("M63").Value = -1487
$ws.Range("N63").Value = -3929.1428

$ws.Range("H66").Value = 2448.795
$ws.Range("I66").Value = 2173
$ws.Range("J66").Value = 2557.1428
$ws.Range("K66").Value = 10865
$ws.Range("L66").Value = 12785.714
$ws.Range("M66").Value = -7433
$ws.Range("N66").Value = -19649.714

$ws.Range("H74").Value = 1123.4814
$ws.Range("I74").Value = 1173.9286
$ws.Range("J74").Value = 1069.1538
$ws.Range("K74").Value = 1173.9286
$ws.Range("L74").Value = 1069.1538
$ws.Range("M74").Value = -299.9286
$ws.Range("N74").Value = -2817.1538

$ws.Range("H77").Value = 1123.4814
$ws.Range("I77").Value = 1173.9286
$ws.Range("J77").Value = 1069.1538
$ws.Range("K77").Value = 5869.643
$ws.Range("L77").Value = 5345.769
$ws.Range("M77").Value = -1501.643
$ws.Range("N77").Value = -14081.769

$ws.Range("H110").Value = 2396.9333
$ws.Range("I110").Value = 1496.7273
$ws.Range("J110").Value = 4872.5
$ws.Range("K110").Value = 1496.7273
$ws.Range("L110").Value = 4872.5
$ws.Range("M110").Value = 548.2727
$ws.Range("N110").Value = -8962.5

$ws.Range("H132").Value = 1201447.8
$ws.Range("I132").Value = 762.8
$ws.Range("J132").Value = 6537825
$ws.Range("K132").Value = 2288.4
$ws.Range("L132").Value = 19613475
$ws.Range("M132").Value = 241.6000000000004
$ws.Range("N132").Value = -19618535

$ws.Range("H136").Value = 9805476
$ws.Range("I136").Value = 9805476
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 29416428
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -29413878
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 12781.833
$ws.Range("I82").Value = 9276.75
$ws.Range("J82").Value = 19792
$ws.Range("K82").Value = 9276.75
$ws.Range("L82").Value = 19792
$ws.Range("M82").Value = -8893.75
$ws.Range("N82").Value = -20558

$ws.Range("H85").Value = 12781.833
$ws.Range("I85").Value = 9276.75
$ws.Range("J85").Value = 19792
$ws.Range("K85").Value = 9276.75
$ws.Range("L85").Value = 19792
$ws.Range("M85").Value = -7950.75
$ws.Range("N85").Value = -22444

$ws.Range("H134").Value = 4635189.5
$ws.Range("I134").Value = 1338.6316
$ws.Range("J134").Value = 22243822
$ws.Range("K134").Value = 4015.8948
$ws.Range("L134").Value = 66731466
$ws.Range("M134").Value = -1480.8948
$ws.Range("N134").Value = -66736536

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1916558.6
$ws.Range("I31").Value = 2058409.2
$ws.Range("J31").Value = 1575
$ws.Range("K31").Value = 2058409.2
$ws.Range("L31").Value = 1575
$ws.Range("M31").Value = -2058114.2
$ws.Range("N31").Value = -2165

$ws.Range("H34").Value = 1916558.6
$ws.Range("I34").Value = 2058409.2
$ws.Range("J34").Value = 1575
$ws.Range("K34").Value = 2058409.2
$ws.Range("L34").Value = 1575
$ws.Range("M34").Value = -2058207.2
$ws.Range("N34").Value = -1979

$ws.Range("H58").Value = 35715028
$ws.Range("I58").Value = 55556320
$ws.Range("J58").Value = 705.9
$ws.Range("K58").Value = 55556320
$ws.Range("L58").Value = 705.9
$ws.Range("M58").Value = -55556117
$ws.Range("N58").Value = -1111.9

$ws.Range("H132").Value = 30306086
$ws.Range("J132").Value = 66670948
$ws.Range("L132").Value = 200012844
$ws.Range("N132").Value = -200017904

$ws.Range("H134").Value = 1358.6111
$ws.Range("I134").Value = 1130.3334
$ws.Range("J134").Value = 2500
$ws.Range("K134").Value = 3391.0002
$ws.Range("L134").Value = 7500
$ws.Range("M134").Value = -856.0001999999999
$ws.Range("N134").Value = -12570

$ws.Range("H136").Value = 35715028
$ws.Range("I136").Value = 55556320
$ws.Range("J136").Value = 705.9
$ws.Range("K136").Value = 166668960
$ws.Range("L136").Value = 2117.7
$ws.Range("M136").Value = -166666410
$ws.Range("N136").Value = -7217.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5905.56
$ws.Range("I132").Value = 1749.5238
$ws.Range("J132").Value = 27724.75
$ws.Range("K132").Value = 5248.5714
$ws.Range("L132").Value = 83174.25
$ws.Range("M132").Value = -2718.5714
$ws.Range("N132").Value = -88234.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 34637636
$ws.Range("I132").Value = 47621010
$ws.Range("J132").Value = 15299.667
$ws.Range("K132").Value = 142863030
$ws.Range("L132").Value = 45899.001
$ws.Range("M132").Value = -142860500
$ws.Range("N132").Value = -50959.001

$ws.Range("H136").Value = 112784290
$ws.Range("I136").Value = 71431160
$ws.Range("K136").Value = 214293480
$ws.Range("M136").Value = -214290930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 44152.223
$ws.Range("I132").Value = 88276
$ws.Range("J132").Value = 8853.200000000001
$ws.Range("K132").Value = 264828
$ws.Range("L132").Value = 26559.6
$ws.Range("M132").Value = -262298
$ws.Range("N132").Value = -31619.6

$ws.Range("H136").Value = 10871542
$ws.Range("I136").Value = 25001250
$ws.Range("J136").Value = 2535.577
$ws.Range("K136").Value = 75003750
$ws.Range("L136").Value = 7606.731000000001
$ws.Range("M136").Value = -75001200
$ws.Range("N136").Value = -12706.731
